$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 306, pushing existing rows 306-369 down to 307-370.
$ws.Rows(306).Insert()

# Populate the newly inserted row 306 with the new record's data.
$ws.Cells.Item(306, 1).Value  = 11
$ws.Cells.Item(306, 2).Value  = "Vega Monumental Concepción"
$ws.Cells.Item(306, 3).Value  = "Bíobío"
$ws.Cells.Item(306, 4).Value  = 44798
$ws.Cells.Item(306, 5).Value  = 8
$ws.Cells.Item(306, 6).Value  = 100112002
$ws.Cells.Item(306, 7).Value  = "Pimiento"
$ws.Cells.Item(306, 8).Value  = "Morrón rojo"
$ws.Cells.Item(306, 9).Value  = "Primera"
$ws.Cells.Item(306, 10).Value = 100
$ws.Cells.Item(306, 11).Value = 25000
$ws.Cells.Item(306, 12).Value = 26000
$ws.Cells.Item(306, 13).Value = 25500
$ws.Cells.Item(306, 14).Value = "$/caja 18 kilos"
$ws.Cells.Item(306, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(306, 16).Value = 1417
$ws.Cells.Item(306, 17).Value = 18
$ws.Cells.Item(306, 18).Value = "Hortaliza"
